$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: extend the "bed" sentence in the Hyponymy/polysemy discussion
# with new text about teaching polysemous words, split across multiple runs
# to match the original diff run-layout.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "involve a flat surface made of a relatively soft material.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "involve a flat surface made of a relatively soft material. When teaching a polysemous word, it may be best to start with the central meaning (or",
    2) | Out-Null

$anchor = $d.Content.Find.Execute(
    "central meaning (or", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null

# Locate the run we just edited so we can append further runs right after it.
$found = $d.Content
$found.Find.Execute("central meaning (or", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($found.End, $found.End)
$insertPoint.InsertAfter(" ")

$found2 = $d.Content
$found2.Find.Execute("central meaning (or ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint2 = $d.Range($found2.End, $found2.End)
$insertPoint2.InsertAfter([char]0x2018)

$found3 = $d.Content
$found3.Find.Execute("central meaning (or " + [char]0x2018, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint3 = $d.Range($found3.End, $found3.End)
$insertPoint3.InsertAfter("sense")

$found4 = $d.Content
$found4.Find.Execute("central meaning (or " + [char]0x2018 + "sense", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint4 = $d.Range($found4.End, $found4.End)
$insertPoint4.InsertAfter([char]0x2019)

$found5 = $d.Content
$found5.Find.Execute("central meaning (or " + [char]0x2018 + "sense" + [char]0x2019, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint5 = $d.Range($found5.End, $found5.End)
$insertPoint5.InsertAfter(") of the word, and later teach the more peripheral meanings (e.g. Tyler & Evans, 2004)")

# ---------------------------------------------------------------------------
# Change 2: split "structure / systematicThe only way..." into two
# sentences, moving "The" to start the next sentence/run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "structure / systematicThe",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "structure / systematic.",
    2) | Out-Null

$d.Content.Find.Execute(
    "only way to test this is via personal intuition.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The only way to test this is via personal intuition.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: add the new Tyler & Evans (2004) bibliography entry as a new
# BodyText paragraph at the end of the document.
# ---------------------------------------------------------------------------
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Last
$lastPara.Style = "BodyText"
$newRange = $lastPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("Tyler, A., & Evans, V. (2004). Applying Cognitive Linguistics to Pedagogical Grammar: The Case of Over.")

$p2 = $d.Paragraphs.Last.Range
$p2.Collapse(0)
$p2.InsertAfter(" ")

$p3 = $d.Paragraphs.Last.Range
$p3.Collapse(0)
$p3.InsertAfter("Cognitive Linguistics, Second Language Acquisition, and Foreign Language Teaching")
$p3.Font.Italic = $true

$p4 = $d.Paragraphs.Last.Range
$p4.Collapse(0)
$p4.InsertAfter(", 257.")
